$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '29.948.87'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.113.40'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '346.12'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").Value = '0.5179'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.4441'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '53.66'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("D10").Value = '0.09381'
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = '1.181'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = '25.21'
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").Value = '8.517'
$ws.Range("E13").Value = '  +4.11%  '
$ws.Range("D14").Value = '2.108.17'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '6.906'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '102.91'
$ws.Range("E16").Value = '  +3.05%  '
$ws.Range("D17").Value = '0.00001163'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '1.009'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = '21.50'
$ws.Range("E19").Value = '  +3.35%  '
$ws.Range("D20").Value = '0.06703'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '6.292'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").Value = '1.007'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '29.986.76'
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '12.71'
$ws.Range("D25").Value = '2.319'
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").Value = '2.371.71'
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("D27").Value = '22.07'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").Value = '2.536'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").Value = '162.58'
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("D30").Value = '134.10'
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").Value = '1.149'
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("D32").Value = '1.781'
$ws.Range("E32").Value = '  +8.12%  '
$ws.Range("D33").Value = '0.1055'
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("D34").Value = '6.237'
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").Value = '6.586'
$ws.Range("E35").Value = '  +5.74%  '
$ws.Range("D36").Value = '3.972'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").Value = '10.80'
$ws.Range("E37").Value = '  +6.03%  '
$ws.Range("D38").Value = '0.02612'
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("D39").Value = '0.06846'
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("D40").Value = '0.7067'
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("D42").Value = '1.333'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '0.2239'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = '0.6829'
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("D45").Value = '14.54'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").Value = '2.354'
$ws.Range("E46").Value = '  +2.29%  '
$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '1.275'
$ws.Range("E48").Value = '  +9.86%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '3.625'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").Value = '0.00000000355'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '1.223'
$ws.Range("E51").Value = '  -0.17%  '
